# Correcion a Diebold Mariano y revision de Cap1
# Insert a new "d=6" row between the existing "d=5" and "d=7" rows,
# pushing the former "d=7" and "d=10" rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move old row 8 ("d=10") down to row 9, preserving the label's formatting.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = $ws.Range("A8").Value()
$ws.Range("B9").Value = $ws.Range("B8").Value()
$ws.Range("C9").Value = $ws.Range("C8").Value()
$ws.Range("D9").Value = $ws.Range("D8").Value()
$ws.Range("E9").Value = $ws.Range("E8").Value()

# Move old row 7 ("d=7") down to row 8.
$ws.Range("A8").Value = $ws.Range("A7").Value()
$ws.Range("B8").Value = $ws.Range("B7").Value()
$ws.Range("C8").Value = $ws.Range("C7").Value()
$ws.Range("D8").Value = $ws.Range("D7").Value()
$ws.Range("E8").Value = $ws.Range("E7").Value()

# Populate row 7 with the new "d=6" entry and its values.
$ws.Range("A7").Value = "d=6"
$ws.Range("B7").Value = 97.85559865433497
$ws.Range("C7").Value = 97.93966211328623
$ws.Range("D7").Value = 97.93993746312258
$ws.Range("E7").Value = 97.94459353283987
